# Build site at 2023-04-12 14:53:07 UTC
# Applies the LOQ4003 sheet content restructuring described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Long text blocks (kept as PowerShell here-strings to preserve embedded
# newlines / ampersands exactly as in the target).
# ---------------------------------------------------------------------------

$objetivosTxt = "Introduzir o aluno na engenharia das reações químicas, através dos conceitos fundamentais da cinética química aplicada a reatores químicos."

$programaResumidoTxt = "Introdução a cinética. Estequiometria cinética. Reações a volume constante. Reações a volume variável. Coleta e análise de dados cinéticos. Métodos de análise e ajuste dos dados cinéticos. Cinética das reações complexas."

$programaTxt = @"
INTRODUÇÃO A CINÉTICA  Tipos de Reações Químicas. Lei de velocidade. Ordem e molecularidade. Constante da velocidade. Tempo de meia-vida e tempo infinito. Influência da temperatura sobre a taxa da reação. Ativação das reações químicas Equação de Arrhenius. Energia de ativação. (4 horas)
ESTEQUIOMETRIA CINÉTICA - Conversão. Concentração e sua variação numa transformação química. (4 horas)
REAÇÕES A VOLUME CONSTANTE: Reações irreversíveis de ordem um. Reações irreversíveis de ordem dois. Reações irreversíveis de ordem três. Reações irreversíveis de ordem qualquer. Reações reversíveis de primeira e segunda ordem. A dependência da constante de equilíbrio com a temperatura. (12 horas)
REAÇÕES A VOLUME VARIÁVEL:  Conceitos. Fração de conversão volumétrica. Reações a volume variável de ordem um e dois. (4 horas)
COLETA E ANÁLISE DE DADOS CINÉTICOS: Introdução. Balanço de massa e coleta de dados em reatores ideais isotérmicos: batelada (BSTR), reator tanque de mistura contínuo (CSTR) e Reator tubular (PFR) (8 horas)
MÉTODOS DE ANÁLISE E AJUSTE DOS DADOS CINÉTICOS: 
Métodos diferencial e integral para o BSTR. Método para o CSTR. Métodos para PFR diferencial e integral. Método das taxas iniciais. Método da meia vida. Método da pressão total (12 horas)
CINÉTICA DAS REAÇÕES COMPLEXAS : Introdução. Mecanismos de reação.  A aproximação do estado estacionário (princípio de Bodenstein). A etapa determinante da velocidade da reação. Reações em cadeia em fase gasosa (Radicais, Pirólise de compostos orgânicos (mecanismo de Rice-Herzfeld), Inibidores e iniciadores, Reações em cadeia ramificada) (8 horas)
CATÁLISE:  Introdução. Catálise homogênea. Catálise heterogênea. Cinética das reações catalíticas heterogêneas. (8 horas)
"@

$criterioTxt = @"
Duas provas escritas (P1 e P2) e trabalhos relacionados à disciplina (TRAB).
"@

$normaRecuperacaoTxt = "Será a média aritmética da nota do aluno na primeira avaliação e da nota do aluno numa prova escrita na recuperação"

$bibliografiaTxt = @"
FOGLER, H. S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2002.
LEVENSPIEL, O. Chemical Reaction Engineering. 3rd. ed. New York: John Wiley & Sons, 1998.
HILL, C.G. An Introduction to chemical engineering kinetics and reactor design. New York: John Wiley&Sons, 1977.
SMITH, J.M. Chemical engineering kinetics. 3rd. ed New York: McGraw-Hill,1981.
DENBIGH, K. ; TURNER, R. Introduction to chemical Reaction Design. Cambridge: Cambridge University Press, 1970.
FROMENT, G.F. ; BISCHOFF, K.B. Chemical reactor analysis and design. 2nd. Ed. New York: John Wiley & Sons, 1990.
"@

$requisitosValTxt = @"
LOQ4088 -  Termodinâmica Química Aplicada II  (Requisito fraco)

"@

# ---------------------------------------------------------------------------
# 1) Update row 10 (Objetivos) B/C with the new objectives text.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = $objetivosTxt
$ws.Range("C10").Value = $objetivosTxt

# ---------------------------------------------------------------------------
# 2) Clear the cells that move out of their old position (the old layout had
#    every A-column label one row "early" relative to its B/C content from
#    row 13 onward; we rebuild rows 13-23 from scratch and add a new row 24).
#    Clear() (not just Value="") removes the cell entirely, matching rows
#    that in the target XML have no <c> element at all for that column.
# ---------------------------------------------------------------------------
$ws.Range("A13").Clear() | Out-Null
$ws.Range("B15").Clear() | Out-Null
$ws.Range("C15").Clear() | Out-Null
$ws.Range("B18").Clear() | Out-Null
$ws.Range("C18").Clear() | Out-Null
$ws.Range("B23").Clear() | Out-Null
$ws.Range("C23").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 3) Cells that keep style "A" (bold label, style index 1) and only need new
#    text - no format copy required.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("A16").Value = "Programa:"
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("A19").Value = "Método:"
$ws.Range("A20").Value = "Critério:"
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("A22").Value = "Bibliografia:"

# ---------------------------------------------------------------------------
# 4) Cells that already carry the correct B/C (wrap, black / wrap, red)
#    styles in the original sheet - just overwrite their value.
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "6310316 - Liana Alvares Rodrigues"
$ws.Range("C13").Value = "6310316 - Liana Alvares Rodrigues"

$ws.Range("B19").Value = $criterioTxt
$ws.Range("C19").Value = $criterioTxt

$ws.Range("B20").Value = "Média da Primeira Avaliação = (I)  Prova P1=50%; (II)  Prova P2=50% e (III) `nObs: Fica a critério de cada docente a inserção de trabalhos no decorrer do curso, bem como a alteração do peso de cada prova em decorrência dos mesmos."
$ws.Range("C20").Value = "Média da Primeira Avaliação = (I)  Prova P1=50%; (II)  Prova P2=50% e (III) `nObs: Fica a critério de cada docente a inserção de trabalhos no decorrer do curso, bem como a alteração do peso de cada prova em decorrência dos mesmos."

$ws.Range("B21").Value = $normaRecuperacaoTxt
$ws.Range("C21").Value = $normaRecuperacaoTxt

# ---------------------------------------------------------------------------
# 5) Cells that are newly populated (were blank / had no style before) -
#    copy number/format (wrap text, vertical alignment, font colour) from a
#    same-column cell that already has the right style, then set the value.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = $programaResumidoTxt
$ws.Range("C14").Value = $programaResumidoTxt

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = $programaTxt
$ws.Range("C16").Value = $programaTxt

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").Value = $bibliografiaTxt
$ws.Range("C22").Value = $bibliografiaTxt

$ws.Range("A3").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = "Requisitos:"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$ws.Range("B24").Value = $requisitosValTxt
$ws.Range("C24").Value = $requisitosValTxt

# ---------------------------------------------------------------------------
# 6) Row heights: restore default (15, no explicit height) for rows that no
#    longer need a tall custom height, and set explicit heights to match the
#    target layout for the others.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).AutoFit() | Out-Null
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit() | Out-Null
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit() | Out-Null
$ws.Rows.Item(24).RowHeight = 30

# ---------------------------------------------------------------------------
# 7) Column layout fix: the original file has a redundant / overlapping
#    column-width definition (col 1-2 @30.71 immediately overridden by a
#    second col 2 @60.71 entry). Nudge column B's width (it keeps its
#    effective value) so the engine splits the stale range and column A is
#    left as its own single-column definition, matching the cleaned-up XML.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 2).EntireColumn.ColumnWidth = 60.7109375

Write-Host "LOQ4003 sheet content rebuilt (rows 10, 13-24) + row heights updated."
